# Generate Report for Handback
#
# This handback-status report is regenerated periodically; each run
# re-stamps the "Xliff Generate Date" / handoff-handback timestamp cells
# with the current run time, while everything else in the workbook
# (file names, statuses, hyperlinks, tables, etc.) stays the same.
#
# Cells touched by this regeneration:
#   Overview!G2   Latest HO Xliff Generate Date   (root .md file)
#   zh-cn!H2      Correspond Handoff Datetime     (f0797733... row)
#   zh-cn!K2      Correspond Handback DateTime    (f0797733... row)
#   de-de!K2      Correspond Handback DateTime    (f0797733... row)

$wb = $excel.ActiveWorkbook

$overview = $wb.Sheets.Item("Overview")
$zhcn     = $wb.Sheets.Item("zh-cn")
$dede     = $wb.Sheets.Item("de-de")

$overview.Range("G2").Value = "2016-09-07 13:56:48"

$zhcn.Range("H2").Value = "2016-09-07 13:56:36"
$zhcn.Range("K2").Value = "2016-09-07 13:57:36"

$dede.Range("K2").Value = "2016-09-07 13:57:56"
